$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44425
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 12500
$ws.Range("Q2").Value = '$/bandeja 18 kilos granel'
$ws.Range("S2").Value = 694
# Row 3
$ws.Range("D3").Value = 44272
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 9500
$ws.Range("Q3").Value = '$/caja 15 kilos granel'
$ws.Range("S3").Value = 633
$ws.Range("T3").Value = 15
# Row 4
$ws.Range("D4").Value = 44272
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 8000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 8000
$ws.Range("Q4").Value = '$/caja 15 kilos granel'
$ws.Range("S4").Value = 533
# Row 5
$ws.Range("D5").Value = 44363
$ws.Range("N5").Value = 9000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 9500
$ws.Range("Q5").Value = '$/caja 15 kilos empedrada'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 633
$ws.Range("T5").Value = 15
# Row 6
$ws.Range("D6").Value = 44307
$ws.Range("L6").Value = 'Primera'
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("Q6").Value = '$/bandeja 18 kilos granel'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 556
# Row 7
$ws.Range("D7").Value = 44307
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 8000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 8000
$ws.Range("Q7").Value = '$/bandeja 18 kilos granel'
$ws.Range("S7").Value = 444
$ws.Range("T7").Value = 18
# Row 8
$ws.Range("D8").Value = 44299
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 11000
$ws.Range("P8").Value = 10500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Región del Maule'
$ws.Range("S8").Value = 583
$ws.Range("T8").Value = 18
# Row 9
$ws.Range("D9").Value = 44299
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 9000
$ws.Range("O9").Value = 9000
$ws.Range("P9").Value = 9000
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("R9").Value = 'Región del Maule'
$ws.Range("S9").Value = 500
# Row 10
$ws.Range("D10").Value = 44316
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 9000
$ws.Range("P10").Value = 9500
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("S10").Value = 528
# Row 11
$ws.Range("D11").Value = 44358
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 11000
$ws.Range("O11").Value = 12000
$ws.Range("P11").Value = 11500
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("S11").Value = 639
